# Adding 5 search test cases
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Fill the new rows in the same order the original author's shared-string
# table shows (A41, B41, C42, C41, A42, B42, then the Runmode/Results cols)
# so new shared strings land on the same indices as the source edit.
$ws.Range("A41").Value = "TestCase_E40"
$ws.Range("B41").Value = "OPQA-1108"
$ws.Range("C42").Value = "Verify that same post can be added to multiple watchlists"
$ws.Range("C41").Value = "Verify that same patent can be added to multiple watchlists"
$ws.Range("A42").Value = "TestCase_E41"
$ws.Range("B42").Value = "OPQA-1109"
$ws.Range("D41").Value = "Y"
$ws.Range("E41").Value = "PASS"
$ws.Range("D42").Value = "Y"
$ws.Range("E42").Value = "PASS"

# Copy style (borders) from row 40 into the new rows
$ws.Range("A40:E40").Copy()
$ws.Range("A41:E42").PasteSpecial(-4122)

# Update selection to reflect the new active cell
$ws.Range("C41").Select()

# Update workbook window size (best effort; host may not persist this)
$excel.ActiveWindow.Width = 12240
$excel.ActiveWindow.Height = 10125
